# Horarios actualizados Línea 141 - 361
# Update timestamps and schedule rows across the three worksheets.

$wb = $excel.ActiveWorkbook

$oldTs = "02:43:45"
$newTs = "02:56:21"

# ---------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: $newTs"

$ws1.Range("A6").Value = $newTs
$ws1.Range("B6").Value = "03:55"
$ws1.Range("C6").Value = "14_ABASTO"
$ws1.Range("D6").Value = 59
$ws1.Range("E6").Value = "LP1912"

$ws1.Range("A7").Value = $newTs
$ws1.Range("B7").Value = "04:01"
$ws1.Range("C7").Value = "81_EL PELIGRO"
$ws1.Range("D7").Value = 65
$ws1.Range("E7").Value = "LP1912"

$ws1.Range("A8").Value = $newTs
$ws1.Range("B8").Value = "04:31"
$ws1.Range("C8").Value = "215_ALUAR"
$ws1.Range("D8").Value = 95
$ws1.Range("E8").Value = "LP1912"

$ws1.Range("A9").Value = $newTs
$ws1.Range("B9").Value = "04:53"
$ws1.Range("C9").Value = "11_ETCHEVERRY"
$ws1.Range("D9").Value = 117
$ws1.Range("E9").Value = "LP1912"

# ---------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: $newTs"
$ws2.Range("A3").Value = "Total filas: 1"

$ws2.Range("A6").Value = $newTs
$ws2.Range("B6").Value = "04:31"
$ws2.Range("C6").Value = "215_ALUAR"
$ws2.Range("D6").Value = 95
$ws2.Range("E6").Value = "LP1912"

# Remove the now-stale second data row (old row 7) so the sheet has only
# one data row left, matching "Total filas: 1".
$ws2.Rows.Item(7).Delete()

# ---------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: $newTs"
